$d = $word.ActiveDocument

# Panel conforme abbreviations
$d.Content.Find.Execute("Member1", $true, $false, $false, $false, $false, $true, 1, $false, "M1", 2)
$d.Content.Find.Execute("PanelChair", $true, $false, $false, $false, $false, $true, 1, $false, "PC", 2)
$d.Content.Find.Execute("Member2", $true, $false, $false, $false, $false, $true, 1, $false, "M3", 2)

# BET5 process: "Membe" + "3" (typo'd "Member3") -> "M4"
$d.Content.Find.Execute("Membe3", $true, $false, $false, $false, $false, $true, 1, $false, "M4", 2)

# "Member" + "4" ("Member4") -> "M5"
$d.Content.Find.Execute("Member4", $true, $false, $false, $false, $false, $true, 1, $false, "M5", 2)
